$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: add J2, update K2 ---
$ws.Range("J2").Value = "generic"
$ws.Range("K2").Value = "can"

# --- Row 3: add J3 ---
$ws.Range("J3").Value = "generic"

# --- Row 4: add J4 ---
$ws.Range("J4").Value = "generic"

# --- Row 5: add J5 ---
$ws.Range("J5").Value = "generic"

# --- Row 6: insert generic into C6, shift D6 stays "can" (same text, but index changed only) ---
$ws.Range("C6").Value = "generic"
$ws.Range("D6").Value = "can"

# --- Row 7 ---
$ws.Range("C7").Value = "generic"
$ws.Range("D7").Value = "can"

# --- Row 8 ---
$ws.Range("C8").Value = "generic"
$ws.Range("D8").Value = "do"

# --- Row 9 ---
$ws.Range("C9").Value = "generic"
$ws.Range("D9").Value = "do"

# --- Row 10 ---
$ws.Range("C10").Value = "generic"
$ws.Range("D10").Value = "look"

# --- Row 11 ---
$ws.Range("C11").Value = "generic"
$ws.Range("D11").Value = "look"

# --- Row 12 ---
$ws.Range("C12").Value = "generic"
$ws.Range("D12").Value = "where"

# --- Row 13 ---
$ws.Range("C13").Value = "generic"
$ws.Range("D13").Value = "where"

# --- New block starting at row 27 ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "video"

$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "video"

$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "video"

$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "video"

$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "audio"

$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "audio"

$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "audio"

$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "audio"
